$d = $word.ActiveDocument

# Locate the paragraph that currently holds the single run of text we need
# to re-split, and collapse a range to its very end (right after the
# trailing "implemented." text, before the paragraph mark).
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "Auto layout constraints for each page were implemented.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)

# Build the replacement paragraph as a WordprocessingML fragment: the same
# sentence, but split across several runs (mirroring how Word splits runs
# across separate edits), with the "_GoBack" bookmark sitting between the
# newly added sentence and the trailing "Auto layout..." sentence, and with
# " except for the images." swapped out for the new
# ". Adding a founder functionality." text. The paragraph keeps the exact
# identity attributes of the original paragraph so it takes its place.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="445C3718" w14:textId="789551FA" w:rsidR="00C3701C" w:rsidRDefault="00F14B06" w:rsidP="00BB1C10"><w:r><w:t>Directory sorts by last name with section headers. You can drag on side to jump to other sections. All data</w:t></w:r><w:r><w:t xml:space="preserve"> saves</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Adding a founder functionality.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> Auto layout constraints for each page were implemented.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$anchor.InsertXML($xml)

# Remove the original (now duplicated) paragraph, including its paragraph
# mark and its copy of the "_GoBack" bookmark, leaving only the newly
# inserted, correctly-split paragraph behind.
$old = $d.Content
$oldFound = $old.Find.Execute(
    "Directory sorts by last name with section headers. You can drag on side to jump to other sections. All data saves except for the images. Auto layout constraints for each page were implemented.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$oldPara = $old.Paragraphs(1)
$oldPara.Range.Delete()
